$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (Synonyms) to make room for "Subspecies".
# This shifts the old F..J headers (Synonyms, Herbarium, Spcode, Liana, binomial)
# one column to the right, to G..K.
$ws.Columns.Item(6).Insert()

# Fill in the new "Subspecies" header in the freshly inserted column F.
$ws.Range("F1").Value = "Subspecies"

# Add the new trailing "fullname" header in column L (one past the shifted "binomial" in K).
$ws.Range("L1").Value = "fullname"

# Remove the two data rows (Swartzia simplex var. grandiflora / var. continentalis),
# leaving only the header row.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
